$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$statQuery = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nWITH f,a,ct,c`n  WHERE a.pubmed_id IN ['31765263']`nRETURN`n    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n    COUNT(DISTINCT c.case_id) AS Cases,`n    COUNT(DISTINCT f) AS Files"

$casesQuery = "MATCH (c:case)`nOPTIONAL  MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)`n MATCH (f:file)-[*]->(c)`n WHERE c.disease = ""Adenocarcinoma of the cervix""`nRETURN DISTINCT`n    c.case_id AS ``Case ID``,`n     ct.clinical_trial_designation AS ``Trial Code``,`n     a.arm_id AS Arm,`n      a.arm_drug AS ``Arm Treatment``,`nc.disease AS Diagnosis,`n  c.gender AS Gender,`n    c.race AS Race,`n    c.ethnicity AS Ethnicity"

$filesQuery = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nOPTIONAL MATCH (f)-->(parent)`nWITH f,a,ct,c,parent`nWHERE c.disease = ""Adenocarcinoma of the cervix""`nWITH`n    f, parent, c, a, ct,`n    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n    toInteger(floor(log(f.file_size)/log(1024))) as i,`n    2 as precision`nWITH`n    f, parent, c, a, ct,`n    f.file_size /(1024^i) AS value,`n    10^precision AS factor,`n    units[i] as unit`nWITH`n    f, parent, c, a, ct, unit,`n    round(factor * value)/factor AS size`nRETURN DISTINCT`n    f.file_name AS ``File Name``,`n    head(labels(parent)) as Association,`n    f.file_description AS Description,`n    f.file_format AS ``File Format``,`n    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    ct.clinical_trial_designation AS ``Trial Code``,`n    a.arm_id AS Arm,`n    c.case_id AS ``Case ID``"

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("B3").Select()
